# "Generate Report for handoff"
#
# The 9eb1fb6a-... file and the 3aaff79e-... file swap table rows (the
# localization-status report re-sorts), and the 3aaff79e-... file's status
# moves from "Handed back: in sync with en-US" to "Ready for handoff" with a
# fresh "Latest Handoff Datetime" stamp, on every sheet (Overview, zh-cn,
# de-de).  Hyperlink relationship targets (r:id -> external URL) are left
# exactly as-is; only the visible display text on each hyperlink moves with
# its cell.

$wb = $excel.ActiveWorkbook

function Swap-Hyperlink {
    param($ws, [string]$cellRef, [string]$target, [string]$display)
    $ws.Range($cellRef).Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $display)
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B3").Value2 = "Ready for handoff"
$ov.Range("C3").Value2 = "Ready for handoff"

Swap-Hyperlink $ov "A2" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$ov.Range("A2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"

Swap-Hyperlink $ov "A3" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md" "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$ov.Range("A3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B3").Value2 = "Ready for handoff"
$zh.Range("D3").Value2 = "2016-02-19 05:46:43"

Swap-Hyperlink $zh "A2" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$zh.Range("A2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"

Swap-Hyperlink $zh "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d74910d327150fa34d2b892b174ffb1eb73e82eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"
$zh.Range("C2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"

Swap-Hyperlink $zh "E2" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/737ebed51987189460a8a5d3ed743f4a02187f9b/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$zh.Range("E2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"

Swap-Hyperlink $zh "F2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/877e69e8d9d0ecbf428f39138039b8f280c728d7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"
$zh.Range("F2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"

Swap-Hyperlink $zh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md" "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$zh.Range("A3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"

Swap-Hyperlink $zh "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d74910d327150fa34d2b892b174ffb1eb73e82eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf" "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"
$zh.Range("C3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"

Swap-Hyperlink $zh "E3" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/737ebed51987189460a8a5d3ed743f4a02187f9b/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md" "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$zh.Range("E3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"

Swap-Hyperlink $zh "F3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/877e69e8d9d0ecbf428f39138039b8f280c728d7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf" "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"
$zh.Range("F3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B3").Value2 = "Ready for handoff"
$de.Range("D3").Value2 = "2016-02-19 05:46:53"

Swap-Hyperlink $de "A2" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$de.Range("A2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"

Swap-Hyperlink $de "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a43525ba630a81de43cbcf7977460be8cb2f356/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"
$de.Range("C2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"

Swap-Hyperlink $de "E2" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/38ca9685ee8aa00cda42d513174a72ed4a4419c6/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$de.Range("E2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"

Swap-Hyperlink $de "F2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/65e9f98ed5eaf37c683148f15b7ada6a021354d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf" "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"
$de.Range("F2").Value2 = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"

Swap-Hyperlink $de "A3" "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md" "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$de.Range("A3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"

Swap-Hyperlink $de "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a43525ba630a81de43cbcf7977460be8cb2f356/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf" "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
$de.Range("C3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"

Swap-Hyperlink $de "E3" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/38ca9685ee8aa00cda42d513174a72ed4a4419c6/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md" "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$de.Range("E3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"

Swap-Hyperlink $de "F3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/65e9f98ed5eaf37c683148f15b7ada6a021354d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf" "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
$de.Range("F3").Value2 = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
